$d = $word.ActiveDocument

# --- Education block -------------------------------------------------
# Collapse the two paragraphs
#   "{{ edu.degree }} - {{ edu.institution }}"
#   "{{ edu.year }}{% if edu.grade %} | {{ edu.grade }}{% endif %}"
# into a single paragraph:
#   "{{ edu.degree }} from {{ edu.institution }} ({{ edu.year }})"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*{{ edu.degree }}*{{ edu.institution }}*") {
        $p.Range.Text = "{{ edu.degree }} from {{ edu.institution }} ({{ edu.year }})"
        $d.Paragraphs.Item($i + 1).Range.Delete()
        break
    }
}

# --- Experience block --------------------------------------------------
# Collapse the two paragraphs
#   "{{ exp.role }} at {{ exp.company }}"
#   "{{ exp.duration }}"
# into a single paragraph:
#   "{{ exp.role }} - {{ exp.company }} ({{ exp.duration }})"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*{{ exp.role }}*{{ exp.company }}*") {
        $p.Range.Text = "{{ exp.role }} - {{ exp.company }} ({{ exp.duration }})"
        $d.Paragraphs.Item($i + 1).Range.Delete()
        break
    }
}

# --- Projects block ------------------------------------------------------
# "Technologies: " -> "Tech: "
$null = $d.Content.Find.Execute("Technologies: ", $true, $false, $false, $false, $false,
                                 $true, 1, $false, "Tech: ", 2)
